# "ajuste em exibicao de data" - append new apostas (bet) rows 27-33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-09-08 13:05:29", 22, 39, 48, 50, 54, 56),
    @("2025-09-08 13:05:30",  5,  6, 14, 22, 44, 48),
    @("2025-09-08 13:05:31",  3, 10, 23, 43, 52, 55),
    @("2025-09-08 13:05:31",  6,  8,  9, 13, 31, 52),
    @("2025-09-08 13:07:41",  4,  6, 35, 36, 43, 55),
    @("2025-09-08 13:07:42",  5, 16, 18, 32, 40, 44),
    @("2025-09-08 13:07:44",  1, 22, 26, 29, 33, 49)
)

$startRow = 27
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A: timestamp stored as literal text, matching the rest of the sheet.
    $ws.Cells.Item($r, 1).Value = [string]$data[0]

    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
